$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "SOIC-8"
$ws.Range("B2").Select()
